$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row values ----
$ws.Range("A1").Value = 'Question'
$ws.Range("B1").Value = 'GroundTruth'
$ws.Range("C1").Value = 'Answer'
$ws.Range("D1").Value = 'Faithfulness_Score'
$ws.Range("E1").Value = 'Faithfulness_Explanation'
$ws.Range("F1").Value = 'Relevance_Score'
$ws.Range("G1").Value = 'Relevance_Explanation'
$ws.Range("H1").Value = 'Context_Precision_Score'
$ws.Range("I1").Value = 'Context_Precision_Explanation'
$ws.Range("J1").Value = 'Context_Recall_Score'
$ws.Range("K1").Value = 'Context_Recall_Explanation'
$ws.Range("L1").Value = 'Context_Documents'

# ---- Row 2 ----
$ws.Range("A2").Value = 'Ẩm thực vùng Cần Giờ có đặc điểm gì nổi bật so với các khu vực khác của Thành phố Hồ Chí Minh?'
$ws.Range("B2").Value = 'Ẩm thực Cần Giờ nổi bật với các món hải sản tươi sống do địa phương có đường bờ biển dài và nguồn thủy hải sản phong phú.'
$ws.Range("C2").Value = 'Ẩm thực Huyện Cần Giờ nổi bật với hải sản tươi ngon và món bóng cá sấy khô chiên.'
$ws.Range("D2").Value = 0.92
$ws.Range("E2").Value = 'The answer states that Cần Giờ''s cuisine is notable for fresh seafood and fried dried fish balls, which directly matches the context that mentions "hải sản và bóng cá sấy khô chiên" as characteristic dishes of Huyện Cần Giờ. The added adjective "tươi ngon" (fresh and tasty) is a reasonable elaboration and does not contradict the source. Therefore the answer is highly faithful to the provided information.'
$ws.Range("F2").Value = 0.88
$ws.Range("G2").Value = 'The answer directly addresses the question by highlighting distinctive culinary features of Cần Giờ—fresh seafood and a specific dish (fried dried fish balls). While brief, it is on-topic and provides the requested comparison point, thus highly relevant.'
$ws.Range("H2").Value = 0.62
$ws.Range("I2").Value = 'The first retrieved document contains the exact statement that Cần Giờ district is known for seafood and fried dried fish balls, directly supporting the answer. However, the remaining documents are unrelated to Cần Giờ’s cuisine, adding noise and lowering overall precision.'
$ws.Range("J2").Value = 0.55
$ws.Range("K2").Value = 'The context includes a line stating that Cần Giờ is known for "hải sản" (seafood) and dried fried fish, which captures the core idea that the district''s cuisine features seafood. However, it does not mention the reasons highlighted in the ground truth—namely the long coastline and abundant marine resources, nor does it emphasize the freshness of the seafood. Therefore, the context partially recalls the key fact but misses important details, resulting in a moderate recall score.'
$ws.Range("L2").Value = 'Tên tài liệu: Ngoài ra, còn có các món đặc trưng của các huyện ngoại thành như thịt bê và gà tre hấp lá giang tại
Nội dung: Ngoài ra, còn có các món đặc trưng của các huyện ngoại thành như thịt bê và gà tre hấp lá giang tại Huyện Củ Chi, hải sản và bóng cá sấy khô chiên ở Huyện Cần Giờ, tóp mỡ ở Huyện Hóc Môn. Huyện Bình Chánh có ẩm thực pha trộn giữa miền Nam và miền Tây. Huyện Nhà Bè giữ được nét ẩm thực miền Nam với cá lóc kho tộ. Quận Bình Tân nổi tiếng với món bún riêu mực. Quận Bình Thạnh có món cháo vịt nổi tiếng ở Thanh Đa. Quận Phú Nhuận có món bánh flan Phan Đăng Lưu.
Score: 0.5542292
Tên tài liệu: Ẩm thực Thành phố Hồ Chí Minh như thế nào?
Nội dung: Ẩm thực Thành phố Hồ Chí Minh như thế nào? Ẩm thực Thành phố Hồ Chí Minh là một bức tranh rực rỡ được dệt nên từ vô số ảnh hưởng văn hóa. Với điều kiện đất đai đa dạng và khí hậu nóng ẩm quanh năm, thành phố có nguồn thực phẩm tươi sống dồi dào. Bức tranh ẩm thực nơi đây phản ánh lịch sử phong phú, khi đã hấp thu hương vị từ khắp Việt Nam và cả thế giới, bao gồm ẩm thực Trung Hoa, Khmer, Thái và châu Âu. Sự giao thoa của những truyền thống ẩm thực này đã tạo nên vô số món ăn hấp dẫn, tiêu biểu cho sự hòa trộn văn hóa của thành phố. Dù là thưởng thức ẩm thực đường phố trong những con hẻm nhộn nhịp, nếm đặc sản tại các quán vỉa hè, hay trải nghiệm ẩm thực cao cấp tại nhà hàng sang trọng, ăn uống ở Sài Gòn luôn là một hành trình hương vị khó quên, để lại dấu ấn sâu đậm trong lòng mỗi thực khách. Những món ăn nhất định phải thử ở Thành phố Hồ Chí Minh 2. Cơm tấm – món ăn nổi tiếng của Thành phố Hồ Chí Minh Đặc điểm: Cơm tấm là sự hòa quyện giữa truyền thống và sáng tạo. Thành phần chính là hạt gạo tấm – sản phẩm phụ của quá trình xay xát. Điểm đặc biệt của cơm tấm chính là sự đa dạng của các món ăn kèm: sườn nướng đậm đà như steak kiểu Âu, bì giòn, chả trứng béo ngậy, nấm mèo luộc, tất cả hòa quyện cùng mỡ hành và nước mắm. Địa điểm gợi ý: Cơm Tấm Ba Ghiền: 84 Đặng Văn Ngữ, P.10, Q.Phú Nhuận Cơm Tấm Mộc: 85 Lý Tự Trọng, P.Bến Thành, Q.1 Cơm Tấm Minh Long: 607 Nguyễn Thị Thập, P.Tân Hưng, Q.7 2.2
Score: 0.51252484
Tên tài liệu: Cẩm nang Ẩm thực TP. Hồ Chí Minh Thành phố Hồ Chí Minh – thành phố sôi động và hiện đại nhất Việt Na
Nội dung: Cẩm nang Ẩm thực TP. Hồ Chí Minh Thành phố Hồ Chí Minh – thành phố sôi động và hiện đại nhất Việt Nam, nổi tiếng với nền văn hóa ẩm thực đầy màu sắc, là nơi hội tụ và giao thoa của nhiều phong cách ẩm thực khác nhau: từ món ăn Việt đến Châu Âu, Nhật Bản, Hàn Quốc, Trung Quốc, Ấn Độ... từ ẩm thực đường phố đến nhà hàng sang trọng. Có lẽ đây là lý do ẩm thực TP. Hồ Chí Minh luôn thu hút sự chú ý đặc biệt, bao gồm cả các kênh truyền hình lớn trên thế giới như CNN Travel, CNN, National Geographic, Discovery... Trong tay bạn là Cẩm nang Ẩm thực TP. Hồ Chí Minh - thành phố của sự sống động, hội tụ và kết hợp từ nhiều vùng miền trên cả nước và nhiều nền văn hóa khác trên thế giới. Hãy cùng tham gia Cẩm nang Ẩm thực TP. Hồ Chí Minh: trải nghiệm văn hóa - thưởng thức ẩm thực độc đáo - và tạo ra những kỷ niệm khó quên cho chuyến đi của bạn. Trải nghiệm văn hóa: Những trải nghiệm chân thực giúp du khách hiểu rõ hơn về văn hóa và lối sống của người dân địa phương. Thưởng thức ẩm thực độc đáo: Với sự hài hòa của các nguyên liệu và sự tinh tế trong cách chế biến, du khách sẽ được thưởng thức những món ăn đặc trưng và độc đáo của TP. Hồ Chí Minh. Khám phá du lịch ẩm thực và tìm hiểu lịch sử: Những hoạt động này mang đến cho du khách những trải nghiệm đáng nhớ và khác biệt trong suốt hành trình của mình.
Score: 0.39963588
Tên tài liệu: Ẩm thực thế giới trong lòng Thành phố Hồ Chí Minh Là nơi hội tụ của cư dân từ mọi miền đất nước cùng
Nội dung: Ẩm thực thế giới trong lòng Thành phố Hồ Chí Minh Là nơi hội tụ của cư dân từ mọi miền đất nước cùng như cửa ngõ tiếp xúc với thế giới bên ngoài nên Thành phố Hồ Chí Minh đã tiếp nhận rất nhiều dòng ẩm thực của cả nước và thế giới, chọn lọc mọi tinh hoa, tạo thành một nền ẩm thực ở Sài Gòn – Thành phố Hồ Chí Minh rất phong phú, đa dạng, hấp dẫn. Tại Sài Gòn – Thành phố Hồ Chí Minh, ngoài những món ngon đến từ các vùng miền khắp Việt Nam, khách du lịch có thể thưởng thức vô số đặc sản trên thế giới với hương vị bản địa ngay tại Sài Gòn. Nói về món ăn Ý thì phải nhắc đến Pizza, mì Ý với rất nhiều loại nước xốt khác nhau xốt cà chua, xốt tôm mực, xốt sữa... Món ăn Đức thì không thể thiếu thịt nguội, jambon, xúc xích vô cùng phong phú và đa dạng. Đặc biệt nhất là bia Đức, được sản xuất ngay tại các nhà hàng. Món ăn Pháp thường ăn kèm xà lách trộn và bánh mì và không thể thiếu rượu vang. Món ăn Hàn Quốc đa số là món nướng ăn kèm với rau xanh với rượu gạo, gần giống rượu nếp của Việt Nam. Món ăn Nhật Bản thì có sushi hay sashimi. Món ăn Trung Hoa thì có Dim Sum (món ăn sáng) hay vịt quay Bắc Kinh. Ngoài ra còn có món ăn Thái Lan, món ăn Ấn Độ, món ăn Hala dành cho người Hồi giáo... Hơn thế nữa, ngay tại Thành phố Hồ Chí Minh, du khách có thể tìm đến một số khu vực, được xem như là “chợ” chuyên kinh doanh các loại gia vị cho các món ăn Campuchia, Hàn Quốc, Nhật Bản, Ấn Độ...
Score: 0.24256428
Tên tài liệu: MÓN ĂN NGON Du khách từ khắp nơi trên thế giới đến TP. Hồ Chí Minh để theo đuổi ước mơ và xây dựng s
Nội dung: MÓN ĂN NGON Du khách từ khắp nơi trên thế giới đến TP. Hồ Chí Minh để theo đuổi ước mơ và xây dựng sự nghiệp, các món ăn truyền thống từ nhiều nơi khác nhau đã có một ngôi nhà mới ở TP. Hồ Chí Minh, bổ sung thêm nhiều màu sắc và hương vị vào bức tranh ẩm thực phong phú của thành phố. Từ ẩm thực đường phố đến các món ăn tinh tế và phức tạp, ẩm thực của thành phố đảm bảo du khách có một trải nghiệm đáng nhớ và đầy hương vị. Bánh mì Việt Nam Bánh tráng trộn Bánh xèo Bún bò Huế Xôi Chè Cơm tấm Hủ tiếu Phở Bò kho
Score: 0.14940722'

# ---- Row 3 ----
$ws.Range("A3").Value = 'Tại khu vực trung tâm Thành phố Hồ Chí Minh, du khách có thể trải nghiệm ẩm thực đường phố ở đâu?'
$ws.Range("B3").Value = 'Du khách có thể trải nghiệm ẩm thực đường phố đa dạng tại khu Phố đi bộ Bùi Viện, nơi tập trung nhiều quán ăn và đồ uống dành cho khách du lịch.'
$ws.Range("C3").Value = 'Bạn có thể thưởng thức ẩm thực đường phố tại các địa điểm sau ở trung tâm Thành phố Hồ Chí Minh:
- **Chợ Bến Thành (Quận 1)** – nơi tập trung nhiều món ăn truyền thống và khu ẩm thực đêm sôi động.  
- **Hẻm 76 Hai Bà Trưng (Quận 1)** – hẻm nhỏ bán các món ăn bình dân như bún Thái, bún chả, súp cua với giá phải chăng.  
- **Khu vực Hồ Con Rùa (Quận 3)** – khu tụ tập về đêm với các món vặt như bánh tráng nướng, gỏi cuốn.  
- **Các gánh hàng rong và quầy bán bánh tráng trộn dọc vỉa hè** quanh trung tâm, nơi du khách có thể tự chọn và thưởng thức ngay tại chỗ.  '
$ws.Range("D3").Value = 0.96
$ws.Range("E3").Value = 'The answer lists Chợ Bến Thành, Hẻm 76 Hai Bà Trưng, Hồ Con Rùa, and street‑vendor stalls, all of which are explicitly mentioned in the provided context as street‑food locations in central Saigon. No information is fabricated or contradicted, so the response is highly faithful to the source.'
$ws.Range("F3").Value = 0.96
$ws.Range("G3").Value = 'The answer directly lists specific locations in central Ho Chi Minh City where tourists can experience street food, exactly addressing the question. It is focused, relevant, and provides useful examples.'
$ws.Range("H3").Value = 0.92
$ws.Range("I3").Value = 'The retrieved contexts contain a document that explicitly lists the same street‑food spots mentioned in the answer (Chợ Bến Thành, Hẻm 76 Hai Bà Trưng, Hồ Con Rùa) and another that describes street‑food stalls in the city centre. Although some other documents are only loosely related, the core information needed to support the answer is present and directly relevant, resulting in a high precision score.'
$ws.Range("J3").Value = 0.05
$ws.Range("K3").Value = 'The context discusses street food in Ho Chi Minh City but never mentions Phố đi bộ Bùi Viện or the specific detail that it offers a diverse range of stalls and drinks for tourists. Therefore it covers almost none of the key facts from the ground truth.'
$ws.Range("L3").Value = 'Tên tài liệu: Bánh tráng trộn Nhắc đến món ăn vặt quen thuộc của người Sài thành, không thể không kể đến bánh trán
Nội dung: Bánh tráng trộn Nhắc đến món ăn vặt quen thuộc của người Sài thành, không thể không kể đến bánh tráng trộn. Đến TP.HCM, du khách có thể dễ dàng bắt gặp những gánh hàng rong hay những cửa hàng nhỏ bán món ăn này dọc vỉa hè, thu hút cả dân văn phòng đến học sinh, sinh viên. Bánh tráng trộn được xem là món ăn vặt đa dạng hương vị như: bánh tráng trộn khô bò, bánh tráng bơ, bánh tráng sa tế, bánh tráng trộn Tây Ninh... Mỗi công thức trộn đều có hương vị riêng biệt, khiến bất kỳ ai cũng có thể bị "nghiện". Du khách có thể dễ dàng khám phá ẩm thực đường phố và cảm nhận được nếp sống sinh hoạt của người dân bản địa. Tại Thành phố Hồ Chí Minh, du khách có thể làm điều đó khi đứng ở vỉa hè chờ người bán làm cho mình một bịch bánh tráng trộn độc đáo.
Score: 0.95292586
Tên tài liệu: C. Các Khu Phố Ẩm Thực Hẻm 200 Xóm Chiếu (Q.4): Thiên đường ẩm thực đường phố, nổi tiếng với các món
Nội dung: C. Các Khu Phố Ẩm Thực Hẻm 200 Xóm Chiếu (Q.4): Thiên đường ẩm thực đường phố, nổi tiếng với các món bún hến, phở, ốc... Khu Chợ Lớn (Q.5 & Q.6): Đặc trưng với ẩm thực người Hoa như vịt quay, sủi cảo, chè, mì xào... Đường ẩm thực Hồ Thị Kỷ (Q.10): Đầy ắp đồ ăn hai bên đường, nổi tiếng với bún ốc hải sản và các món cà ri Khmer. Đường Phan Xích Long (Q. Phú Nhuận): Mệnh danh là "thiên đường ẩm thực" đa dạng hương vị Bắc - Trung - Nam, với nhiều nhà hàng nổi tiếng. Hẻm 76 Hai Bà Trưng (Q.1): Hẻm nhỏ bán các món ăn bình dân, giá cả phải chăng như bún Thái, bún chả, súp cua. Hồ Con Rùa (Q.3): Khu vực tụ tập về đêm với các món ăn vặt đường phố giá phải chăng như bánh tráng nướng, gỏi cuốn. Chợ Bến Thành (Q.1) và Chợ Tân Định (Q.1): Cung cấp đa dạng các món ăn truyền thống và khu ẩm thực đêm sôi động.
Score: 0.87610537
Tên tài liệu: Ẩm thực Thành phố Hồ Chí Minh như thế nào?
Nội dung: Ẩm thực Thành phố Hồ Chí Minh như thế nào? Ẩm thực Thành phố Hồ Chí Minh là một bức tranh rực rỡ được dệt nên từ vô số ảnh hưởng văn hóa. Với điều kiện đất đai đa dạng và khí hậu nóng ẩm quanh năm, thành phố có nguồn thực phẩm tươi sống dồi dào. Bức tranh ẩm thực nơi đây phản ánh lịch sử phong phú, khi đã hấp thu hương vị từ khắp Việt Nam và cả thế giới, bao gồm ẩm thực Trung Hoa, Khmer, Thái và châu Âu. Sự giao thoa của những truyền thống ẩm thực này đã tạo nên vô số món ăn hấp dẫn, tiêu biểu cho sự hòa trộn văn hóa của thành phố. Dù là thưởng thức ẩm thực đường phố trong những con hẻm nhộn nhịp, nếm đặc sản tại các quán vỉa hè, hay trải nghiệm ẩm thực cao cấp tại nhà hàng sang trọng, ăn uống ở Sài Gòn luôn là một hành trình hương vị khó quên, để lại dấu ấn sâu đậm trong lòng mỗi thực khách. Những món ăn nhất định phải thử ở Thành phố Hồ Chí Minh 2. Cơm tấm – món ăn nổi tiếng của Thành phố Hồ Chí Minh Đặc điểm: Cơm tấm là sự hòa quyện giữa truyền thống và sáng tạo. Thành phần chính là hạt gạo tấm – sản phẩm phụ của quá trình xay xát. Điểm đặc biệt của cơm tấm chính là sự đa dạng của các món ăn kèm: sườn nướng đậm đà như steak kiểu Âu, bì giòn, chả trứng béo ngậy, nấm mèo luộc, tất cả hòa quyện cùng mỡ hành và nước mắm. Địa điểm gợi ý: Cơm Tấm Ba Ghiền: 84 Đặng Văn Ngữ, P.10, Q.Phú Nhuận Cơm Tấm Mộc: 85 Lý Tự Trọng, P.Bến Thành, Q.1 Cơm Tấm Minh Long: 607 Nguyễn Thị Thập, P.Tân Hưng, Q.7 2.2
Score: 0.86759764
Tên tài liệu: Ẩm thực thế giới trong lòng Thành phố Hồ Chí Minh Là nơi hội tụ của cư dân từ mọi miền đất nước cùng
Nội dung: Ẩm thực thế giới trong lòng Thành phố Hồ Chí Minh Là nơi hội tụ của cư dân từ mọi miền đất nước cùng như cửa ngõ tiếp xúc với thế giới bên ngoài nên Thành phố Hồ Chí Minh đã tiếp nhận rất nhiều dòng ẩm thực của cả nước và thế giới, chọn lọc mọi tinh hoa, tạo thành một nền ẩm thực ở Sài Gòn – Thành phố Hồ Chí Minh rất phong phú, đa dạng, hấp dẫn. Tại Sài Gòn – Thành phố Hồ Chí Minh, ngoài những món ngon đến từ các vùng miền khắp Việt Nam, khách du lịch có thể thưởng thức vô số đặc sản trên thế giới với hương vị bản địa ngay tại Sài Gòn. Nói về món ăn Ý thì phải nhắc đến Pizza, mì Ý với rất nhiều loại nước xốt khác nhau xốt cà chua, xốt tôm mực, xốt sữa... Món ăn Đức thì không thể thiếu thịt nguội, jambon, xúc xích vô cùng phong phú và đa dạng. Đặc biệt nhất là bia Đức, được sản xuất ngay tại các nhà hàng. Món ăn Pháp thường ăn kèm xà lách trộn và bánh mì và không thể thiếu rượu vang. Món ăn Hàn Quốc đa số là món nướng ăn kèm với rau xanh với rượu gạo, gần giống rượu nếp của Việt Nam. Món ăn Nhật Bản thì có sushi hay sashimi. Món ăn Trung Hoa thì có Dim Sum (món ăn sáng) hay vịt quay Bắc Kinh. Ngoài ra còn có món ăn Thái Lan, món ăn Ấn Độ, món ăn Hala dành cho người Hồi giáo... Hơn thế nữa, ngay tại Thành phố Hồ Chí Minh, du khách có thể tìm đến một số khu vực, được xem như là “chợ” chuyên kinh doanh các loại gia vị cho các món ăn Campuchia, Hàn Quốc, Nhật Bản, Ấn Độ...
Score: 0.85269034
Tên tài liệu: MÓN ĂN NGON Du khách từ khắp nơi trên thế giới đến TP. Hồ Chí Minh để theo đuổi ước mơ và xây dựng s
Nội dung: MÓN ĂN NGON Du khách từ khắp nơi trên thế giới đến TP. Hồ Chí Minh để theo đuổi ước mơ và xây dựng sự nghiệp, các món ăn truyền thống từ nhiều nơi khác nhau đã có một ngôi nhà mới ở TP. Hồ Chí Minh, bổ sung thêm nhiều màu sắc và hương vị vào bức tranh ẩm thực phong phú của thành phố. Từ ẩm thực đường phố đến các món ăn tinh tế và phức tạp, ẩm thực của thành phố đảm bảo du khách có một trải nghiệm đáng nhớ và đầy hương vị. Bánh mì Việt Nam Bánh tráng trộn Bánh xèo Bún bò Huế Xôi Chè Cơm tấm Hủ tiếu Phở Bò kho
Score: 0.84619683'

# ---- Row 4 ----
$ws.Range("A4").Value = 'Chợ Bình Tây ở Quận 6 nổi tiếng với các loại hàng hóa và ẩm thực như thế nào?'
$ws.Range("B4").Value = 'Chợ Bình Tây là nơi buôn bán sầm uất, có nhiều gian hàng ẩm thực, hàng khô, trái cây và các món ăn đặc trưng của người Hoa.'
$ws.Range("C4").Value = 'Chợ Bình Tây ở Quận 6 nổi tiếng với một dải sản phẩm đa dạng: thực phẩm khô, gia vị, bánh kẹo, mứt, các mặt hàng nhu yếu phẩm, quần áo, giày dép, đồ gia dụng và đồ lưu niệm. Ngoài các sạp bán buôn, chợ còn có khu ẩm thực phong phú với các xe ẩm thực hoạt động vào buổi tối, các món ăn truyền thống và các lễ hội, sự kiện văn hóa thu hút du khách trong và ngoài nước.'
$ws.Range("D4").Value = 0.99
$ws.Range("E4").Value = 'The answer accurately lists the product categories (dry food, spices, sweets, jams, necessities, clothing, shoes, household goods, souvenirs) and the night‑time food‑truck area, festivals and cultural events that attract tourists, all of which are explicitly mentioned in the provided context. No unsupported or fabricated information is present.'
$ws.Range("F4").Value = 0.96
$ws.Range("G4").Value = 'The answer directly addresses the question by listing the types of goods (dry food, spices, sweets, clothing, etc.) and describing the market''s culinary offerings (food stalls, night trucks, traditional dishes, festivals). It is focused and provides a comprehensive overview, making it highly relevant.'
$ws.Range("H4").Value = 0.98
$ws.Range("I4").Value = 'The retrieved passages directly mention the wide range of goods sold at Chợ Bình Tây—including spices, sweets, jams, clothing, shoes, household items, and souvenirs—as well as the vibrant night‑time food stalls, festivals and cultural events that attract tourists. This information closely matches the answer, making the context highly precise.'
$ws.Range("J4").Value = 0.78
$ws.Range("K4").Value = 'The context mentions that Chợ Bình Tây is a large, historic market with thousands of stalls, a wide variety of food items, dry goods, and a strong Chinese community, covering most of the ground truth points (busy market, many food stalls, dry goods). However, it does not explicitly mention fruit or specifically highlight the characteristic Chinese dishes, and the term "sầm uất" (bustling) is not directly used. Therefore, the recall is high but not complete.'
$ws.Range("L4").Value = 'Tên tài liệu: CHỢ BÌNH TÂY: Với khuôn viên rộng 000m$^2$, nằm giữa tuyến đường Tháp Mười - Lê Tấn Kế - Phan Văn Kh
Nội dung: CHỢ BÌNH TÂY: Với khuôn viên rộng 000m$^2$, nằm giữa tuyến đường Tháp Mười - Lê Tấn Kế - Phan Văn Khỏe - Trần Bình, Chợ Bình Tây mở ra một thế giới đa dạng về ẩm thực, thực phẩm, đồ lưu niệm và kể cả những sản phẩm truyền thống khó tìm thấy ở nơi khác nhưng chỉ với mức giá phải chăng. Bên cạnh lối kiến trúc cổ xưa cùng bề dày lịch sử, hiện nay Chợ Bình Tây đã mở rộng và phát triển thêm các hoạt động kinh tế - văn hóa về đêm như tổ chức cho các xe ẩm thực hoạt động kinh doanh, diễn ra một số lễ hội, sự kiện văn hóa tại khu vực trước sân chợ để phục vụ cho du khách trong và ngoài nước. Hằng năm, chợ thu hút trên 1000 lượt khách du lịch nước ngoài đến tham quan và mua sắm tại chợ. Được xây dựng vào những năm 20 của thế kỷ XX, Chợ Bình Tây tự hào là biểu tượng kiến trúc lịch sử của Sài Gòn xưa. Không những vậy, dù đã trải qua thời gian dài hình thành và phát triển, Chợ Bình Tây ngày nay vẫn giữ được vị thế của một trong những chợ đầu mối buôn bán lớn của thành phố và của quận 6. Chợ Bình Tây không chỉ khơi dậy sự tò mò của du khách nước ngoài mà còn đánh thức lòng tự hào của người dân Việt Nam khi lưu giữ những nét văn hóa của phố phường xưa và kiến trúc cổ kính, đậm chất thơ. Địa chỉ: 57A Tháp Mười, phường 2, quận 6.
Score: 0.9936365
Tên tài liệu: CHỢ BÌNH TÂY Với gần 100 năm lịch sử, chợ Bình Tây đại diện cho lịch sử lâu đời và văn hóa độc đáo c
Nội dung: CHỢ BÌNH TÂY Với gần 100 năm lịch sử, chợ Bình Tây đại diện cho lịch sử lâu đời và văn hóa độc đáo của cộng đồng người Hoa. Chợ được chia thành 5 khu vực, với tối đa 2.300 quầy hàng cùng hơn 30 mặt hàng khác nhau, chủ yếu là gia vị, bánh kẹo, mứt, quần áo, giày dép, v.v. 57A Tháp Mười, Phường 2, Quận 6
Score: 0.9792749
Tên tài liệu: Chợ Bình Tây (Chợ Lớn Mới)
Nội dung: . Năm 1975, chính quyền Thành phố tiếp nhận quản lý, sắp xếp cho nhân dân tiếp tục mua bán phục vụ hàng hoá cho cả nước và các nước bạn Lào, Campuchia, Trung Quốc và đổi tên chợ là chợ Bình Tây cho đến ngày nay. Năm 1992, tiếp tục phát huy thế mạnh của chợ, Ủy ban nhân dân Quận 6 tổ chức sửa chữa nâng cấp nhà lồng chợ thêm một tầng lầu. Năm 2006, tiếp tục đầu tư cải tạo sửa chữa khu vực đường Trần Bình - Lê Tấn Kế khang trang sạch đẹp, chợ Bình Tây vì thế trở thành một trong những ngôi chợ lớn của Thành phố với 2.358 quầy sạp. Khu vực nhà lồng chợ có 446 sạp, trong đó tầng trệt là 698 sạp, tầng lầu có 748 sạp. Khu vực ngoài nhà lồng có 912 sạp, trong đó khu vực đường Trần Bình có 408 sạp, khu vực đường Lê Tấn Kế có 328 sạp, khu vực đường Phan Văn Khoẻ là 176 sạp. Các ngành hàng được bố trí, sắp xếp hợp lý, tập trung theo từng khu vực kinh doanh nhằm phát huy thế mạnh của từng ngành hàng. Bà con người Hoa vẫn tập trung về chợ làm ăn mua bán mà đa số là người Hoa sinh sống tại các Quận 5, Quận 6 và Quận 1 Tiểu thương người Hoa hiện chiếm tỷ lệ 25% số lượng hộ kinh doanh tại chợ Bình Tây. Trải qua thời gian dài hình thành và phát triển, chợ Bình Tây ngày nay vẫn giữ được vị thế của một chợ đầu mối bán buôn lớn của Thành phố và của Quận 6, với lối kiến trúc cổ xưa của Trung Quốc và bề dày lịch sử, hàng năm thu hút trên 1000 lượt khách du lịch người nước ngoài đến tham quan và mua sắm tại chợ
Score: 0.96820134
Tên tài liệu: Top Những Chợ Lớn Nhất Ở TP. Hồ Chí Minh TP. Hồ Chí Minh được xem là một trong những trung tâm kinh tế, văn hóa, du lịch hàng đầu của Việt Nam.
Nội dung: . 2. Chợ Bình Tây (Chợ Lớn) Địa chỉ: 57A Tháp Mười, Phường 2, Quận 6, TP. Hồ Chí Minh Giờ mở cửa: 5:00 – 19:30 Được xây dựng năm 1928 nhờ sự tài trợ của thương nhân người Hoa Quách Đàm, chợ Bình Tây mang đậm kiến trúc Á Đông đặc trưng. Với diện tích 000 m² và hơn 2.300 sạp hàng, đây là trung tâm bán sỉ lớn nhất Sài Gòn với các mặt hàng: thực phẩm khô, nhu yếu phẩm, quần áo, giày dép, đồ gia dụng. 3. Chợ An Đông Địa chỉ: 36 An Dương Vương, Phường 9, Quận 5, TP. Hồ Chí Minh Giờ mở cửa: 7:00 – 18:00 Ra đời năm 1954, An Đông được xem là “thủ phủ thời trang” của Sài Gòn, cùng với Bến Thành và Bình Tây là ba khu chợ lâu đời nhất thành phố. Với hơn 2.000 sạp trên diện tích 000 m², nơi đây tập trung quần áo, phụ kiện, vải vóc nhập từ Trung Quốc, Nhật Bản, Thái Lan… Ngoài ra còn có khu ẩm thực và hàng khô nổi tiếng. 4. Chợ Bà Chiểu Địa chỉ: 40 Điện Biên Phủ, Phường 1, Quận Bình Thạnh, TP. Hồ Chí Minh Giờ mở cửa: 5:00 – 22:00 Có từ thế kỷ 19, được xây dựng lại năm 1942. Với 800 sạp hàng trên diện tích 8.465 m², chợ Bà Chiểu nổi tiếng nhất là chợ đồ si lớn nhất Sài Gòn, nơi có nhiều sản phẩm từ các thương hiệu quốc tế với giá rẻ. Về đêm, chợ Bà Chiểu còn là “thiên đường ẩm thực” với cơm gà, bún cá, bún thịt nướng, chè 6000 đồng… 5. Chợ Tân Định Địa chỉ: 336 Hai Bà Trưng, Phường Tân Định, Quận 1, TP. Hồ Chí Minh Giờ mở cửa: 5:00 – 18:00 Được xây dựng năm 1927, chợ Tân Định từng nổi tiếng là “chợ nhà giàu” do giá cả cao hơn so với nơi khác
Score: 0.92455155
Tên tài liệu: Top Những Chợ Lớn Nhất Ở TP. Hồ Chí Minh TP. Hồ Chí Minh được xem là một trong những trung tâm kinh tế, văn hóa, du lịch hàng đầu của Việt Nam.
Nội dung: Top Những Chợ Lớn Nhất Ở TP. Hồ Chí Minh TP. Hồ Chí Minh được xem là một trong những trung tâm kinh tế, văn hóa, du lịch hàng đầu của Việt Nam. Thành phố này không chỉ nổi bật với tốc độ phát triển kinh tế mạnh mẽ mà còn gây ấn tượng với sự đa dạng về ẩm thực và văn hóa. Một trong những điểm thu hút khách du lịch chính là các khu chợ truyền thống. Các khu chợ này bày bán đủ loại sản phẩm, từ thực phẩm, quần áo, mỹ phẩm, đồ gia dụng cho đến các mặt hàng đặc sản. Dưới đây là danh sách 9 khu chợ lớn nhất TP. Hồ Chí Minh, theo gợi ý của Visithcmc. Chợ Bến Thành Địa chỉ: Phường Bến Thành, Quận 1, TP. Hồ Chí Minh Giờ mở cửa: 7:30 – 18:00 Chợ Bến Thành là biểu tượng nổi tiếng nhất của thành phố, đặc biệt thu hút khách quốc tế. Được xây dựng từ năm 1912, chợ có diện tích 13.056 m² với 4 cổng chính và 12 cổng phụ mở ra bốn con đường sầm uất. Bên trong có khoảng 500 gian hàng, bán đủ loại mặt hàng: nhu yếu phẩm, quần áo, trang sức, mỹ phẩm, thủ công mỹ nghệ… Đặc biệt là khu ẩm thực với cơm tấm, bún riêu, bún mắm, bún thịt nướng, chả giò và nhiều món tráng miệng. Ngoài mua sắm, cổng chính chợ Bến Thành còn là “điểm check-in” quen thuộc của du khách. 2. Chợ Bình Tây (Chợ Lớn) Địa chỉ: 57A Tháp Mười, Phường 2, Quận 6, TP. Hồ Chí Minh Giờ mở cửa: 5:00 – 19:30 Được xây dựng năm 1928 nhờ sự tài trợ của thương nhân người Hoa Quách Đàm, chợ Bình Tây mang đậm kiến trúc Á Đông đặc trưng
Score: 0.79486305'

# ---- Row 5 ----
$ws.Range("A5").Value = 'Du khách có thể tìm thấy khu vực bán ẩm thực phong phú nào trong Chợ Bến Thành, Quận 1?'
$ws.Range("B5").Value = 'Trong Chợ Bến Thành có khu ẩm thực với nhiều món ăn truyền thống Việt Nam và món ăn đặc trưng của Sài Gòn phục vụ du khách trong và ngoài nước.'
$ws.Range("C5").Value = 'Trong Chợ Bến Thành, du khách sẽ tìm thấy “khu ẩm thực” – một khu vực tập trung đa dạng các món ăn đặc sản của thành phố như cơm tấm, bún riêu, bún mắm, bún thịt nướng, chả giò và nhiều món tráng miệng. Đây là nơi được xem là thiên đường ẩm thực của quận 1.'
$ws.Range("D5").Value = 0.99
$ws.Range("E5").Value = 'The answer correctly states that Ben Thanh Market has a "khu ẩm thực" (food area) featuring dishes such as cơm tấm, bún riêu, bún mắm, bún thịt nướng, chả giò and many desserts, and that it is considered a food paradise of District 1. These details are directly supported by the provided context passages.'
$ws.Range("F5").Value = 0.7
$ws.Range("G5").Value = 'The answer correctly identifies that there is a "khu ẩm thực" (food area) within Ben Thanh Market where tourists can find a variety of dishes, directly addressing the question. However, it does not provide a more specific location or name beyond the generic term, so the relevance is good but not fully detailed.'
$ws.Range("H5").Value = 0.97
$ws.Range("I5").Value = 'The retrieved passages explicitly describe a dedicated "khu ẩm thực" (food area) inside Chợ Bến Thành, listing the same dishes mentioned in the answer (cơm tấm, bún riêu, bún mắm, bún thịt nướng, chả giò, etc.) and call the market a "thiên đường ẩm thực của quận 1." This directly supports the answer, making the context highly precise.'
$ws.Range("J5").Value = 0.98
$ws.Range("K5").Value = 'The context clearly states that Chợ Bến Thành is a culinary hub in District 1, offering a wide variety of traditional Vietnamese dishes and Saigon specialties that serve both local and foreign tourists, directly matching the ground‑truth information.'
$ws.Range("L5").Value = 'Tên tài liệu: CHỢ BẾN THÀNH: Trải qua thăng trầm của lịch sử và nhiều lần tu sửa, nét cổ kính của Chợ Bến Thành vẫ
Nội dung: CHỢ BẾN THÀNH: Trải qua thăng trầm của lịch sử và nhiều lần tu sửa, nét cổ kính của Chợ Bến Thành vẫn được lưu giữ và trở thành một biểu tượng không thể tách rời khi nhắc đến Thành phố Hồ Chí Minh. Chợ Bến Thành hiện nay có tổng diện tích 13.056m$^2$, chiều dài từ cửa Nam đến cửa Bắc 136m, từ cửa Đông đến cửa Tây 96m, lối đi nội các của Đông, Tây, Nam, Bắc rộng 5m. Bốn cửa lớn đặt tại bốn hướng đều có các phù điêu gồm thể hiện các sản vật địa phương: Cửa Nam - của chính của chợ có hình cá đuối, cá trê, bò và heo; cửa Bắc có hình chuỗi, ngỗng; cửa Đông có hình bò, heo; và cửa Tây có hình cá đuối, chuỗi. Ngoài bốn cửa chính còn có nhiều cửa phụ để ra vào chợ. Nổi bật trên nóc của chính cửa chợ là tháp đồng hồ ba mặt - một nét biểu trưng giúp du khách có thể dễ dàng nhận biết Chợ Bến Thành khi nhìn từ xa. Hiện nay chợ tập trung khoảng gần 500 hộ kinh doanh buôn bán nhiều mặt hàng thiết yếu và độc lạ trong nước cũng như nước ngoài. Nơi đây còn được xem là thiên đường ẩm thực của quận 1, nơi tập trung hầu hết các món ăn đặc sản của thành phố với hương vị thơm ngon, ấn tượng như bún riêu, bún mắm, bánh canh cua, bánh bèo... hoặc các món ăn vật trứ danh của người Sài thành như gỏi cuốn, ốc, bột chiên, xôi 7 màu, chè... Địa chỉ: Quảng trường Quách Thị Trang, phường Bến Thành, quận
Score: 0.9851576
Tên tài liệu: C. Các Khu Phố Ẩm Thực Hẻm 200 Xóm Chiếu (Q.4): Thiên đường ẩm thực đường phố, nổi tiếng với các món
Nội dung: C. Các Khu Phố Ẩm Thực Hẻm 200 Xóm Chiếu (Q.4): Thiên đường ẩm thực đường phố, nổi tiếng với các món bún hến, phở, ốc... Khu Chợ Lớn (Q.5 & Q.6): Đặc trưng với ẩm thực người Hoa như vịt quay, sủi cảo, chè, mì xào... Đường ẩm thực Hồ Thị Kỷ (Q.10): Đầy ắp đồ ăn hai bên đường, nổi tiếng với bún ốc hải sản và các món cà ri Khmer. Đường Phan Xích Long (Q. Phú Nhuận): Mệnh danh là "thiên đường ẩm thực" đa dạng hương vị Bắc - Trung - Nam, với nhiều nhà hàng nổi tiếng. Hẻm 76 Hai Bà Trưng (Q.1): Hẻm nhỏ bán các món ăn bình dân, giá cả phải chăng như bún Thái, bún chả, súp cua. Hồ Con Rùa (Q.3): Khu vực tụ tập về đêm với các món ăn vặt đường phố giá phải chăng như bánh tráng nướng, gỏi cuốn. Chợ Bến Thành (Q.1) và Chợ Tân Định (Q.1): Cung cấp đa dạng các món ăn truyền thống và khu ẩm thực đêm sôi động.
Score: 0.97313875
Tên tài liệu: Ẩm thực Quận 1 Quận 1 là nơi giao thoa của nhiều nền văn hóa ẩm thực, từ châu Âu như Ý, Tây Ban Nha,
Nội dung: Ẩm thực Quận 1 Quận 1 là nơi giao thoa của nhiều nền văn hóa ẩm thực, từ châu Âu như Ý, Tây Ban Nha, Pháp cho đến ẩm thực Nhật Bản. Một số nhà hàng nổi tiếng ở đây là: Din Ky Restaurant: 137B-137C Nguyễn Trãi, P. Bến Thành, Q.1, TP. HCM The Running Bean: 33 Mạc Thị Bưởi, P. Bến Nghé, Q.1, TP. HCM Ann Quan: 66 Ngô Đức Kế, P. Bến Nghé, Q.1, TP. HCM L''Etoile Restaurant: 180 Bis Hai Bà Trưng, P. Đa Kao, Q.1, TP. HCM
Score: 0.7932662
Tên tài liệu: Top Những Chợ Lớn Nhất Ở TP. Hồ Chí Minh TP. Hồ Chí Minh được xem là một trong những trung tâm kinh tế, văn hóa, du lịch hàng đầu của Việt Nam.
Nội dung: Top Những Chợ Lớn Nhất Ở TP. Hồ Chí Minh TP. Hồ Chí Minh được xem là một trong những trung tâm kinh tế, văn hóa, du lịch hàng đầu của Việt Nam. Thành phố này không chỉ nổi bật với tốc độ phát triển kinh tế mạnh mẽ mà còn gây ấn tượng với sự đa dạng về ẩm thực và văn hóa. Một trong những điểm thu hút khách du lịch chính là các khu chợ truyền thống. Các khu chợ này bày bán đủ loại sản phẩm, từ thực phẩm, quần áo, mỹ phẩm, đồ gia dụng cho đến các mặt hàng đặc sản. Dưới đây là danh sách 9 khu chợ lớn nhất TP. Hồ Chí Minh, theo gợi ý của Visithcmc. Chợ Bến Thành Địa chỉ: Phường Bến Thành, Quận 1, TP. Hồ Chí Minh Giờ mở cửa: 7:30 – 18:00 Chợ Bến Thành là biểu tượng nổi tiếng nhất của thành phố, đặc biệt thu hút khách quốc tế. Được xây dựng từ năm 1912, chợ có diện tích 13.056 m² với 4 cổng chính và 12 cổng phụ mở ra bốn con đường sầm uất. Bên trong có khoảng 500 gian hàng, bán đủ loại mặt hàng: nhu yếu phẩm, quần áo, trang sức, mỹ phẩm, thủ công mỹ nghệ… Đặc biệt là khu ẩm thực với cơm tấm, bún riêu, bún mắm, bún thịt nướng, chả giò và nhiều món tráng miệng. Ngoài mua sắm, cổng chính chợ Bến Thành còn là “điểm check-in” quen thuộc của du khách. 2. Chợ Bình Tây (Chợ Lớn) Địa chỉ: 57A Tháp Mười, Phường 2, Quận 6, TP. Hồ Chí Minh Giờ mở cửa: 5:00 – 19:30 Được xây dựng năm 1928 nhờ sự tài trợ của thương nhân người Hoa Quách Đàm, chợ Bình Tây mang đậm kiến trúc Á Đông đặc trưng
Score: 0.73392886
Tên tài liệu: Ẩm thực thế giới trong lòng Thành phố Hồ Chí Minh Là nơi hội tụ của cư dân từ mọi miền đất nước cùng
Nội dung: Ẩm thực thế giới trong lòng Thành phố Hồ Chí Minh Là nơi hội tụ của cư dân từ mọi miền đất nước cùng như cửa ngõ tiếp xúc với thế giới bên ngoài nên Thành phố Hồ Chí Minh đã tiếp nhận rất nhiều dòng ẩm thực của cả nước và thế giới, chọn lọc mọi tinh hoa, tạo thành một nền ẩm thực ở Sài Gòn – Thành phố Hồ Chí Minh rất phong phú, đa dạng, hấp dẫn. Tại Sài Gòn – Thành phố Hồ Chí Minh, ngoài những món ngon đến từ các vùng miền khắp Việt Nam, khách du lịch có thể thưởng thức vô số đặc sản trên thế giới với hương vị bản địa ngay tại Sài Gòn. Nói về món ăn Ý thì phải nhắc đến Pizza, mì Ý với rất nhiều loại nước xốt khác nhau xốt cà chua, xốt tôm mực, xốt sữa... Món ăn Đức thì không thể thiếu thịt nguội, jambon, xúc xích vô cùng phong phú và đa dạng. Đặc biệt nhất là bia Đức, được sản xuất ngay tại các nhà hàng. Món ăn Pháp thường ăn kèm xà lách trộn và bánh mì và không thể thiếu rượu vang. Món ăn Hàn Quốc đa số là món nướng ăn kèm với rau xanh với rượu gạo, gần giống rượu nếp của Việt Nam. Món ăn Nhật Bản thì có sushi hay sashimi. Món ăn Trung Hoa thì có Dim Sum (món ăn sáng) hay vịt quay Bắc Kinh. Ngoài ra còn có món ăn Thái Lan, món ăn Ấn Độ, món ăn Hala dành cho người Hồi giáo... Hơn thế nữa, ngay tại Thành phố Hồ Chí Minh, du khách có thể tìm đến một số khu vực, được xem như là “chợ” chuyên kinh doanh các loại gia vị cho các món ăn Campuchia, Hàn Quốc, Nhật Bản, Ấn Độ...
Score: 0.3192275'

# ---- Reset auto row-height bookkeeping from the multi-line Context_Documents
# text back to the sheet default (matches the source file, which carries no
# custom row heights even though several cells wrap many lines). ----
$ws.Rows(2).AutoFit()
$ws.Rows(3).AutoFit()
$ws.Rows(4).AutoFit()
$ws.Rows(5).AutoFit()
# ---- Header style: bold font + thin box border + center/top alignment ----
# Build the combined format on A1 alone (keeps the style table minimal,
# like a single new cellXf), then fan it out with a formats-only paste so
# every header cell shares that one xf instead of each accreting its own.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)
